$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-08 07:18:20'
$ws.Range('O2').Value = '-2.6 °C'
$ws.Range('E3').Value = '2026-02-08 07:18:22'
$ws.Range('E4').Value = '2026-02-08 07:18:24'
$ws.Range('H4').Value = '74%'
$ws.Range('J4').Value = '1001.5 hPa'
$ws.Range('O4').Value = '7.9 °C'
$ws.Range('E5').Value = '2026-02-08 07:18:27'
$ws.Range('O5').Value = '-4.5 °C'
$ws.Range('E6').Value = '2026-02-08 07:18:29'
$ws.Range('H6').Value = '70%'
$ws.Range('N6').Value = '5.4 °C 6:49 TU'
$ws.Range('O6').Value = '8.1 °C'
$ws.Range('E7').Value = '2026-02-08 07:18:32'
$ws.Range('H7').Value = '72%'
$ws.Range('J7').Value = '1001.4 hPa'
$ws.Range('N7').Value = '9.9 °C 6:45 TU'
$ws.Range('O7').Value = '11.0 °C'
$ws.Range('E8').Value = '2026-02-08 07:18:34'
$ws.Range('J8').Value = '1001.5 hPa'
$ws.Range('L8').Value = '36.7 km/h - 303º 6:59 TU'
$ws.Range('O8').Value = '7.9 °C'
$ws.Range('E9').Value = '2026-02-08 07:18:37'
$ws.Range('H9').Value = '77%'
$ws.Range('O9').Value = '7.3 °C'
$ws.Range('E10').Value = '2026-02-08 07:18:39'
$ws.Range('N10').Value = '3.5 °C 6:58 TU'
$ws.Range('O10').Value = '6.6 °C'
$ws.Range('E11').Value = '2026-02-08 07:18:41'
$ws.Range('N11').Value = '-0.2 °C 6:57 TU'
$ws.Range('E12').Value = '2026-02-08 07:18:44'
$ws.Range('H12').Value = '78%'
$ws.Range('N12').Value = '6.8 °C 6:59 TU'
$ws.Range('O12').Value = '8.5 °C'
$ws.Range('E13').Value = '2026-02-08 07:18:46'
$ws.Range('H13').Value = '93%'
$ws.Range('N13').Value = '-1.3 °C 6:58 TU'
$ws.Range('O13').Value = '0.3 °C'
$ws.Range('E14').Value = '2026-02-08 07:18:48'
$ws.Range('N14').Value = '5.8 °C 6:59 TU'
$ws.Range('O14').Value = '8.0 °C'
$ws.Range('E15').Value = '2026-02-08 07:18:51'
$ws.Range('O15').Value = '6.3 °C'
$ws.Range('E16').Value = '2026-02-08 07:18:53'
$ws.Range('H16').Value = '84%'
$ws.Range('E17').Value = '2026-02-08 07:18:56'
$ws.Range('N17').Value = '-1.5 °C 6:59 TU'
$ws.Range('O17').Value = '-0.9 °C'
$ws.Range('E18').Value = '2026-02-08 07:18:58'
$ws.Range('H18').Value = '80%'
$ws.Range('J18').Value = '1001.6 hPa'
$ws.Range('N18').Value = '4.6 °C 6:59 TU'
$ws.Range('O18').Value = '7.5 °C'
$ws.Range('E19').Value = '2026-02-08 07:19:00'
$ws.Range('M19').Value = '2.8 °C 6:49 TU'
$ws.Range('E20').Value = '2026-02-08 07:19:03'
$ws.Range('N20').Value = '-5.5 °C 6:30 TU'
$ws.Range('O20').Value = '-5.0 °C'
$ws.Range('E21').Value = '2026-02-08 07:19:05'
$ws.Range('E22').Value = '2026-02-08 07:19:07'
$ws.Range('N22').Value = '-8.3 °C 6:59 TU'
$ws.Range('O22').Value = '-6.8 °C'
$ws.Range('E23').Value = '2026-02-08 07:19:10'
$ws.Range('E24').Value = '2026-02-08 07:19:12'
$ws.Range('H24').Value = '91%'
$ws.Range('J24').Value = '1002.2 hPa'
$ws.Range('N24').Value = '2.9 °C 6:59 TU'
$ws.Range('O24').Value = '6.2 °C'
$ws.Range('E25').Value = '2026-02-08 07:19:15'
$ws.Range('H25').Value = '85%'
$ws.Range('E26').Value = '2026-02-08 07:19:17'
$ws.Range('H26').Value = '77%'
$ws.Range('N26').Value = '-0.1 °C 6:59 TU'
$ws.Range('O26').Value = '1.4 °C'
$ws.Range('E27').Value = '2026-02-08 07:19:20'
$ws.Range('H27').Value = '91%'
$ws.Range('E28').Value = '2026-02-08 07:19:22'
$ws.Range('N28').Value = '3.0 °C 6:51 TU'
$ws.Range('O28').Value = '5.3 °C'
$ws.Range('E29').Value = '2026-02-08 07:19:25'
$ws.Range('H29').Value = '93%'
$ws.Range('N29').Value = '5.9 °C 6:51 TU'
$ws.Range('O29').Value = '8.8 °C'
$ws.Range('E30').Value = '2026-02-08 07:19:27'
$ws.Range('H30').Value = '68%'
$ws.Range('N30').Value = '6.7 °C 6:49 TU'
$ws.Range('O30').Value = '8.9 °C'
$ws.Range('E31').Value = '2026-02-08 07:19:30'
$ws.Range('H31').Value = '61%'
$ws.Range('N31').Value = '9.0 °C 6:34 TU'
$ws.Range('O31').Value = '9.9 °C'
$ws.Range('E32').Value = '2026-02-08 07:19:32'
$ws.Range('H32').Value = '100%'
$ws.Range('E33').Value = '2026-02-08 07:19:34'
$ws.Range('H33').Value = '92%'
$ws.Range('N33').Value = '-1.3 °C 6:48 TU'
$ws.Range('O33').Value = '0.3 °C'
$ws.Range('E34').Value = '2026-02-08 07:19:37'
$ws.Range('E35').Value = '2026-02-08 07:19:39'
$ws.Range('J35').Value = '1002.6 hPa'
$ws.Range('O35').Value = '3.3 °C'
$ws.Range('E36').Value = '2026-02-08 07:19:42'
$ws.Range('H36').Value = '73%'
$ws.Range('J36').Value = '1001.4 hPa'
$ws.Range('N36').Value = '8.6 °C 6:59 TU'
$ws.Range('O36').Value = '10.5 °C'
$ws.Range('E37').Value = '2026-02-08 07:19:44'
$ws.Range('N37').Value = '1.3 °C 6:58 TU'
$ws.Range('O37').Value = '3.0 °C'
$ws.Range('E38').Value = '2026-02-08 07:19:47'
$ws.Range('H38').Value = '85%'
$ws.Range('N38').Value = '4.2 °C 6:59 TU'
$ws.Range('O38').Value = '7.1 °C'
$ws.Range('E39').Value = '2026-02-08 07:19:49'
$ws.Range('E40').Value = '2026-02-08 07:19:51'
$ws.Range('J40').Value = '1004.1 hPa'
$ws.Range('E41').Value = '2026-02-08 07:19:54'
$ws.Range('J41').Value = '1001.2 hPa'
$ws.Range('O41').Value = '8.5 °C'
$ws.Range('E42').Value = '2026-02-08 07:19:56'
$ws.Range('N42').Value = '5.5 °C 6:55 TU'
$ws.Range('O42').Value = '8.8 °C'
$ws.Range('E43').Value = '2026-02-08 07:19:58'
$ws.Range('O43').Value = '4.1 °C'
$ws.Range('E44').Value = '2026-02-08 07:20:01'
$ws.Range('E45').Value = '2026-02-08 07:20:03'
$ws.Range('H45').Value = '69%'
$ws.Range('J45').Value = '1002.1 hPa'
$ws.Range('O45').Value = '2.5 °C'
$ws.Range('E46').Value = '2026-02-08 07:20:06'
$ws.Range('H46').Value = '88%'
$ws.Range('O46').Value = '6.1 °C'
